$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.055.93"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.957.95"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.02"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.60"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +1.84%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.586"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.46"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0852"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "3.429.82"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("B14").Value = "Uniswap"
$ws.Range("C14").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.45"
$ws.Range("E14").Value = "  +74.94%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.42"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.75"
$ws.Range("E16").Value = "  +5.67%  "
$ws.Range("D17").Value = "2.957.26"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("E18").Value = "  +4.38%  "
$ws.Range("D19").Value = "51.094.55"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("E20").Value = "  -2.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.42"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").Value = "0.0₃0955"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.36"
$ws.Range("E23").Value = "  +17.54%  "
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "267.54"
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.98"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.82"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.02"
$ws.Range("E30").Value = "  -7.43%  "
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.50"
$ws.Range("E32").Value = "  +7.25%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.16"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.77"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.03"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("E36").Value = "  -4.55%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("E38").Value = "  +8.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.78"
$ws.Range("E39").Value = "  +3.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.117"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("E42").Value = "  -2.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.68"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.59"
$ws.Range("E44").Value = "  +11.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.69"
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.02"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").Value = "2.042.32"
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.257"
$ws.Range("E49").Value = "  -4.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0322"
$ws.Range("E50").Value = "  -5.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.35"
$ws.Range("E51").Value = "  +6.92%  "
